# Updated units for Platelets and WBC after PhysioNet forum posting. Improve analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row: new columns D:H ----
$ws.Cells.Item(1,4).Value = "RangeMin"
$ws.Cells.Item(1,5).Value = "RangeMax"
$ws.Cells.Item(1,6).Value = "MinPlotValue"
$ws.Cells.Item(1,7).Value = "MaxPlotValue"
$ws.Cells.Item(1,8).Value = "Source"

$wustl = "http://idgateway.wustl.edu/Normal%20lab%20values.pdf"
$wiki  = "https://en.wikipedia.org/wiki/Fraction_of_inspired_oxygen"

# ---- Row 2: BaseExcess ----
$ws.Cells.Item(2,6).Value = -12
$ws.Cells.Item(2,7).Value = 12

# ---- Row 3: HCO3 ----
$ws.Cells.Item(3,4).Value = 23
$ws.Cells.Item(3,5).Value = 28
$ws.Cells.Item(3,6).Value = 10
$ws.Cells.Item(3,7).Value = 40
$ws.Cells.Item(3,8).Value = $wustl

# ---- Row 4: FiO2 ----
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = $wiki

# ---- Row 5: pH (units cleared) ----
$ws.Range("B5").Value = ""
$ws.Cells.Item(5,4).Value = 7.38
$ws.Cells.Item(5,5).Value = 7.44
$ws.Cells.Item(5,6).Value = 7
$ws.Cells.Item(5,7).Value = 7.7
$ws.Cells.Item(5,8).Value = $wustl

# ---- Row 6: PaCO2 ----
$ws.Cells.Item(6,4).Value = 35
$ws.Cells.Item(6,5).Value = 45
$ws.Cells.Item(6,6).Value = 10
$ws.Cells.Item(6,7).Value = 75
$ws.Cells.Item(6,8).Value = $wustl

# ---- Row 7: SaO2 ----
$ws.Cells.Item(7,4).Value = 80
$ws.Cells.Item(7,5).Value = 100
$ws.Cells.Item(7,6).Value = 90
$ws.Cells.Item(7,7).Value = 100
$ws.Cells.Item(7,8).Value = $wustl

# ---- Row 8: AST ----
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 250

# ---- Row 9: BUN ----
$ws.Cells.Item(9,4).Value = 8
$ws.Cells.Item(9,5).Value = 20
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 50
$ws.Cells.Item(9,8).Value = $wustl

# ---- Row 10: Alkalinephos ----
$ws.Cells.Item(10,4).Value = 36
$ws.Cells.Item(10,5).Value = 92
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 300
$ws.Cells.Item(10,8).Value = $wustl

# ---- Row 11: Calcium ----
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 12

# ---- Row 12: Chloride ----
$ws.Cells.Item(12,6).Value = 85
$ws.Cells.Item(12,7).Value = 125

# ---- Row 13: Creatinine ----
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 5

# ---- Row 14: Bilirubin_direct ----
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 6

# ---- Row 15: Glucose ----
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 300

# ---- Row 16: Lactate ----
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 7.5

# ---- Row 17: Magnesium ----
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 3.5

# ---- Row 18: Phosphate ----
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 10

# ---- Row 19: Potassium ----
$ws.Cells.Item(19,6).Value = 2
$ws.Cells.Item(19,7).Value = 7.5

# ---- Row 20: Bilirubin_total (description gains trailing space) ----
$ws.Range("C20").Value = "Total bilirubin "
$ws.Cells.Item(20,6).Value = 0
$ws.Cells.Item(20,7).Value = 6

# ---- Row 21: TroponinI ----
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 2

# ---- Row 22: Hct ----
$ws.Cells.Item(22,6).Value = 15
$ws.Cells.Item(22,7).Value = 50

# ---- Row 23: Hgb ----
$ws.Cells.Item(23,6).Value = 5
$ws.Cells.Item(23,7).Value = 17.5

# ---- Row 24: PTT ----
$ws.Cells.Item(24,6).Value = 10
$ws.Cells.Item(24,7).Value = 100

# ---- Row 25: WBC (unit change: count/L -> count*10^3/uL) ----
$ws.Range("B25").Value = "count*10^3/µL"
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 30

# ---- Row 26: Fibrinogen ----
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 1000

# ---- Row 27: Platelets (unit change: count/mL -> count*10^3/uL) ----
$ws.Range("B27").Value = "count*10^3/µL"
$ws.Cells.Item(27,6).Value = 0
$ws.Cells.Item(27,7).Value = 600

# ---- Column widths (approximate Excel's auto-fit for the updated layout;
#      the host's ColumnWidth setter quantizes internally, so inputs are
#      pre-compensated to land as close as possible to the authored widths) ----
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 44.166666666666664
$ws.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 53.5

# ---- Selection ends on C10, matching the recorded author selection ----
$ws.Range("C10").Select()
